# accelerometer_selected.xlsx (falling class) — slide the sample window:
# prepend the two newest accelerometer readings and drop the three oldest
# rows that age out of the window (net: one fewer data row overall).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right below the header (new rows 2 and 3), pushing
# all of the existing samples down by two rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Inserting next to the bold/bordered header row makes Excel copy that
# header formatting onto the new rows by default — strip it back off so
# the new samples look like every other plain data row.
$ws.Range("A2:C3").ClearFormats()

# Write in the two newest accelerometer samples (x, y, z).
$ws.Range("A2").Value = -0.2734694480895995
$ws.Range("B2").Value = 0.2277572751045226
$ws.Range("C2").Value = -0.1111783366650344

$ws.Range("A3").Value = -0.2649335861206055
$ws.Range("B3").Value = 0.1057968139648434
$ws.Range("C3").Value = -0.4681921228766454

# Drop the three oldest samples that fell off the end of the window (they
# now sit at rows 22-24 after the insert above shifted everything by two).
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
